$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PSA_LOLO")
$ws.Range("A2").Value = 41297
$ws.Range("B2").Value = 23479
